$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('2025-07-25T15:16:01+00:00', 'EXTERNAL:- BBC Asian Network - Wk30 - 2025-07-31 - Thursday', 'pressportal@bbc.co.uk', 1)
    ,@('2025-07-25T15:16:00+00:00', 'EXTERNAL:- BBC Radio 4 FM - Wk29 - 2025-07-25 - Friday', 'pressportal@bbc.co.uk', 1)
    ,@('2025-07-25T15:15:59+00:00', 'EXTERNAL:- BBC World Service UK Schedule - Wk31 - 2025-08-02 - Saturday', 'pressportal@bbc.co.uk', 1)
    ,@('2025-07-25T15:15:59+00:00', 'EXTERNAL:- Pop Schedule for Week 31 AMENDMENT', 'schedules@globallistings-distribution.com', 1)
    ,@('2025-07-25T15:15:58+00:00', 'EXTERNAL:- BBC Radio Scotland Extra - Wk30 - 2025-07-26 - Saturday', 'pressportal@bbc.co.uk', 1)
    ,@('2025-07-25T15:15:58+00:00', 'EXTERNAL:- BBC Radio 1 - Wk30 - 2025-07-30 - Wednesday', 'pressportal@bbc.co.uk', 1)
    ,@('2025-07-25T15:15:57+00:00', 'EXTERNAL:- BBC Radio Shetland - Wk29 - 2025-07-25 - Friday', 'pressportal@bbc.co.uk', 1)
    ,@('2025-07-25T15:15:56+00:00', 'EXTERNAL:- BBC Radio Scotland - Wk31 - 2025-08-02 - Saturday', 'pressportal@bbc.co.uk', 1)
    ,@('2025-07-25T15:15:56+00:00', 'EXTERNAL:- BBC Radio 4 Extra - Wk30 - 2025-07-28 - Monday', 'pressportal@bbc.co.uk', 1)
    ,@('2025-07-25T15:15:55+00:00', 'EXTERNAL:- BBC Radio Scotland Extra - Wk31 - 2025-08-02 - Saturday', 'pressportal@bbc.co.uk', 1)
    ,@('2025-07-25T15:15:55+00:00', 'EXTERNAL:- BBC Asian Network - Wk30 - 2025-07-29 - Tuesday', 'pressportal@bbc.co.uk', 1)
    ,@('2025-07-25T15:15:55+00:00', 'EXTERNAL:- BBC Radio Orkney - Wk31 - 2025-08-02 - Saturday', 'pressportal@bbc.co.uk', 1)
    ,@('2025-07-25T15:15:55+00:00', 'EXTERNAL:- BBC Radio 4 FM - Wk30 - 2025-08-01 - Friday', 'pressportal@bbc.co.uk', 1)
    ,@('2025-07-25T15:15:54+00:00', 'EXTERNAL:- BBC Radio Orkney - Wk30 - 2025-08-01 - Friday', 'pressportal@bbc.co.uk', 1)
    ,@('2025-07-25T15:15:53+00:00', 'EXTERNAL:- BBC Radio Wales Extra - Wk30 - 2025-08-01 - Friday', 'pressportal@bbc.co.uk', 1)
    ,@('2025-07-25T15:15:52+00:00', 'EXTERNAL:- BBC Radio 1Xtra - Wk30 - 2025-07-27 - Sunday', 'pressportal@bbc.co.uk', 1)
    ,@('2025-07-25T15:15:51+00:00', 'EXTERNAL:- BBC Radio 4 FM - Wk30 - 2025-07-27 - Sunday', 'pressportal@bbc.co.uk', 1)
    ,@('2025-07-25T15:15:50+00:00', 'EXTERNAL:- BBC Radio Scotland - Wk30 - 2025-07-26 - Saturday', 'pressportal@bbc.co.uk', 1)
    ,@('2025-07-25T15:15:49+00:00', 'EXTERNAL:- BBC Radio 4 FM - Wk30 - 2025-07-31 - Thursday', 'pressportal@bbc.co.uk', 1)
    ,@('2025-07-25T15:15:47+00:00', 'EXTERNAL:- BBC Radio 4 FM - Wk30 - 2025-07-30 - Wednesday', 'pressportal@bbc.co.uk', 1)
    ,@('2025-07-25T15:15:47+00:00', 'EXTERNAL:- BBC Asian Network - Wk29 - 2025-07-25 - Friday', 'pressportal@bbc.co.uk', 1)
    ,@('2025-07-25T15:15:46+00:00', 'EXTERNAL:- BBC Asian Network - Wk30 - 2025-07-26 - Saturday', 'pressportal@bbc.co.uk', 1)
    ,@('2025-07-25T15:15:46+00:00', 'EXTERNAL:- BBC Radio Scotland Extra - Wk29 - 2025-07-25 - Friday', 'pressportal@bbc.co.uk', 1)
    ,@('2025-07-25T15:15:45+00:00', 'EXTERNAL:- BBC Radio Scotland - Wk29 - 2025-07-25 - Friday', 'pressportal@bbc.co.uk', 1)
    ,@('2025-07-25T15:15:45+00:00', 'EXTERNAL:- BBC Radio 1 - Wk29 - 2025-07-25 - Friday', 'pressportal@bbc.co.uk', 1)
    ,@('2025-07-25T15:15:43+00:00', 'EXTERNAL:- BBC Radio 1Xtra - Wk29 - 2025-07-25 - Friday', 'pressportal@bbc.co.uk', 1)
    ,@('2025-07-25T15:15:42+00:00', 'EXTERNAL:- BBC Radio York schedule - Week 35', 'nathan.moore@bbc.co.uk', 1)
    ,@('2025-07-25T15:14:15+00:00', 'EXTERNAL:- SCHEDULE AMENDMENTS: Comedy Central & Comedy Central Extra - Week 31 - 32', 'Daniel.McNamara@paramount.com', 1)
    ,@('2025-07-25T15:11:12+00:00', 'U&Eden - Week 33 Listings', 'Patsy.Keating@pa.media', 1)
    ,@('2025-07-25T15:09:43+00:00', 'U&Gold - Week 33 Listings', 'Manisha.Anand@pa.media', 1)
    ,@('2025-07-25T15:09:05+00:00', 'EXTERNAL:- E4 Extra Week 32 (Sun, Fri) Post Listings Update', 'MPoole@Channel4.co.uk', 0)
    ,@('2025-07-25T15:04:14+00:00', 'EXTERNAL:- For Prep -  MONDAY 28TH JULY ', 'Chanelle.Stratford@qvc.com', 1)
    ,@('2025-07-25T15:00:29+00:00', 'EXTERNAL:- EPG Schedules for SONYTVASIACOM [HD] for 01-07-2025 to 31-07-2025', 'absolute@absmail.tv', 1)
)

$startRow = 2532
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = [bool]$rec[3]
}

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D2564"))
